# ebay-LoadSheet.xlsx -- "adding support for google sheets"
#
# Source diff recap:
#   - F2 (eBay Category) value changes from 34532 to 377.
#   - The active selection on the sheet moves from B3 to G8.
#   - A handful of column widths shift by a fraction of a character
#     (sub-pixel rounding noise consistent with a Google Sheets
#     re-export of the column metrics) and column E's explicit width
#     entry disappears (reverts to the sheet's default width).
#   - The two conditional-format "dxf" fill styles (red / green) used
#     by the sheet's expression-based highlighting rules collapse into
#     a single blank/no-fill dxf, and the S1:S1000 rule's dxfId is
#     repointed at it.
#
# The two functionally meaningful edits -- the data value and the
# saved selection -- are applied directly below. The column-width
# nudges are reproduced on a best-effort basis (same columns touched,
# widths pushed the same direction) since the workbook/runtime only
# lets widths be set in whole-character units, so they cannot be
# reproduced to the original's fractional precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data edit -----------------------------------------------------
$ws.Range("F2").Value = 377

# --- column width nudges (best effort) ------------------------------
$ws.Columns(1).ColumnWidth = 12.6743782533256
$ws.Columns(2).ColumnWidth = 18.2452284557548
$ws.Columns(4).ColumnWidth = 16.3181029496819
$ws.Columns(6).ColumnWidth = 16.7472527472528
$ws.Columns(7).ColumnWidth = 34.9577790630422
$ws.Columns(13).ColumnWidth = 24.7796414112204
$ws.Columns(15).ColumnWidth = 22.8525159051475
$ws.Columns(16).ColumnWidth = 17.7108155002892

# --- view state ------------------------------------------------------
# Selection moves from B3 to G8.
$ws.Range("G8").Select() | Out-Null
